# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give H1 the same header style as the other headers (copy style from G1),
# then set its label.
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(1, 8).Value = "Save"

# Populate the Save values for each data row (rows 2-11).
$saveValues = @(0, 0, 1, 0, 0, 1, 0, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

$excel.CutCopyMode = 0
